$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("a new two-story building (maps were created by", $true, $false, $false, $false, $false, $true, 1, $false, "new area with two floors (maps were created by", 2)
